$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo'd date in the "Chart with realtime data" row (row 21, Source column)
$ws.Range("E21").Value = "#experts-BIRT email sent 7/23/2014 2:46 PM Eastern"

# Row 22: Filter and Search from Header
$ws.Range("A22").Value = "Filter and Search from Header"
$ws.Range("B22").Value = "JSAPI"
$ws.Range("C22").Value = "Filter a table based on a search for a string either within the entire table, or within a single column.  Script in onContentUpdate and BeforeFactory"
$ws.Range("D22").Value = "Bill Clark"
$ws.Range("E22").Value = "#experts-BIRT email sent 7/29/2014"
$ws.Range("F22").Value = "Report Designs/JSAPI/FastFilterDetailDemo.rptdesign"

# Row 23: Change report name in title bar
$ws.Range("A23").Value = "Change report  name in title bar"
$ws.Range("B23").Value = "Scripting"
$ws.Range("C23").Value = "Changes the name that is displayed in the title bar through scripting in OnContentUpdate"
$ws.Range("D23").Value = "Rob Murphy"
$ws.Range("E23").Value = "#experts-BIRT email sent 7/15/2014"
$ws.Range("F23").Value = "Report Designs/Scripting/ChangeReportNameInTitleBar.rptdesign"

# Row 24: Control the size and placement of tooltip
$ws.Range("A24").Value = "Control the size and placement of tooltip"
$ws.Range("B24").Value = "Charts"
$ws.Range("C24").Value = "Uses scripting to control the size and placement of the tooltip on an HTML5 chart.  The scripting is in the script tab on the chart."
$ws.Range("D24").Value = "Clement Wong"
$ws.Range("E24").Value = "#experts-BIRT email sent 6/30/2014"
$ws.Range("F24").Value = "Report Designs/Charts/ChartTooltip__cwong.rptdesign"

# Row 25: Drilldown to any state map
$ws.Range("A25").Value = "Drilldown to any state map"
$ws.Range("B25").Value = "Maps"
$ws.Range("C25").Value = "This shows how to drill down from a map of the US to a map of any state.  The DynamicStateMapByCounty report design takes a state short name as a parameter.  In the beforeFactory script, the report design's map property for the XMLRepresentation is updated based on the parameter value."
$ws.Range("D25").Value = "Glenn Hess, Pierre Tessier"
$ws.Range("F25").Value = "Report Designs/Maps/US Populations by State Map.rptdesign"

# Update the view: selection ends on A25 (also clears the old topLeftCell="D1" scroll position)
$ws.Range("A25").Select()
